$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Aselmo-Manut. Mot. End.-1NA, Anderson-Tornearia-1NA, Valmir-Tec. Mat. Não Metal.-1NA, Anderson-Metrologia 1-1NA]"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "[Elcio Dec.-Desenho tecnico mecanico – T2-1NA, Elcio Dec.-Desenho tecnico mecanico – T2-1NA]"
$ws.Range("E18").Value = "[Rachel-Trat. Termicos-1NA, Emerson-Comandos Eletricos-1NA, Victor-Ajustagem-1NA, J. Paulo S.-Tec. Mat. Não Metal.-1NA]"
$ws.Range("F18").Value = "[Gisele-E. D. N. D.-1NA, Suzanny-Metalografia-1NA, Clesidson-Elet. Dig. Bas.-1NA, Joel L.-T. M. Metalicos-1NA]"

$ws.Range("B19").Value = "[Aselmo-Manut. Mot. End.-1NA, Anderson-Tornearia-1NA, Valmir-Tec. Mat. Não Metal.-1NA, Anderson-Metrologia 1-1NA]"
$ws.Range("C19").Value = "Gilberto-M.T.R.M.-"
$ws.Range("D19").Value = "[Elcio Dec.-Desenho tecnico mecanico – T2-1NA, Elcio Dec.-Desenho tecnico mecanico – T2-1NA]"
$ws.Range("E19").Value = "[Rachel-Trat. Termicos-1NA, Emerson-Comandos Eletricos-1NA, Victor-Ajustagem-1NA, J. Paulo S.-Tec. Mat. Não Metal.-1NA]"
$ws.Range("F19").Value = "[Gisele-E. D. N. D.-1NA, Suzanny-Metalografia-1NA, Clesidson-Elet. Dig. Bas.-1NA, Joel L.-T. M. Metalicos-1NA]"

$ws.Range("B20").Value = "[Humberto-Desenho tecnico mecanico – T1-1NA, Humberto-Desenho tecnico mecanico – T1-1NA]"
$ws.Range("C20").Value = "Tiago P.-M.T.F.-"
$ws.Range("D20").Value = "[Valmir-Tec. Mat. Não Metal.-1NA, Anderson-Tornearia-1NA, Aselmo-Manut. Mot. End.-1NA, Anderson-Metrologia 1-1NA]"
$ws.Range("E20").Value = "[Rachel-Trat. Termicos-1NA, Emerson-Comandos Eletricos-1NA, Victor-Ajustagem-1NA, J. Paulo S.-Tec. Mat. Não Metal.-1NA]"
$ws.Range("F20").Value = "[Gisele-E. D. N. D.-1NA, Suzanny-Metalografia-1NA, Clesidson-Elet. Dig. Bas.-1NA, Joel L.-T. M. Metalicos-1NA]"

$ws.Range("B21").Value = "[Humberto-Desenho tecnico mecanico – T1-1NA, Humberto-Desenho tecnico mecanico – T1-1NA]"
$ws.Range("C21").Value = "Tiago P.-M.T.F.-"
$ws.Range("D21").Value = "[Valmir-Tec. Mat. Não Metal.-1NA, Anderson-Tornearia-1NA, Aselmo-Manut. Mot. End.-1NA, Anderson-Metrologia 1-1NA]"
$ws.Range("E21").Value = "[Rachel-Trat. Termicos-1NA, Emerson-Comandos Eletricos-1NA, Victor-Ajustagem-1NA, J. Paulo S.-Tec. Mat. Não Metal.-1NA]"
$ws.Range("F21").Value = "[Gisele-E. D. N. D.-1NA, Suzanny-Metalografia-1NA, Clesidson-Elet. Dig. Bas.-1NA, Joel L.-T. M. Metalicos-1NA]"
